# "Redo naming of mixres units"
#
# Before: area_mixres, area_mixres_new, area_hires, area_lores, area_pop_sum
# After:  area_mixre, area_hires, area_lores, area_pop_sum
#
# The old "area_mixres" sheet (the first sheet) is dropped entirely, and
# "area_mixres_new" (whose stats then become the new first sheet's data)
# is renamed to "area_mixre". Every other sheet keeps its data and simply
# shifts up one position.

$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation prompt Excel normally shows.
$excel.DisplayAlerts = $false | Out-Null

$wb.Worksheets.Item("area_mixres").Delete() | Out-Null
$wb.Worksheets.Item("area_mixres_new").Name = "area_mixre"

$excel.DisplayAlerts = $true | Out-Null
